# Apply the daily crypto-price refresh (GitHub Actions data pull).
# Columns: A=rank(unchanged) B=Coin C=Link D=Price E=Volume(1h)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "63.181.68"
$ws.Range("E2").Value = "  +0.16%  "
# Row 3
$ws.Range("D3").Value = "2.564.90"
$ws.Range("E3").Value = "  +0.91%  "
# Row 4
$ws.Range("E4").Value = "  +0.19%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "583.92"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.67%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.48"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.62%  "
# Row 7
$ws.Range("E7").Value = "  +0.15%  "
# Row 8
$ws.Range("E8").Value = "  +0.82%  "
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.109"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.01%  "
# Row 10
$ws.Range("E10").Value = "  -0.79%  "
# Row 11
$ws.Range("E11").Value = "  +0.30%  "
# Row 12
$ws.Range("E12").Value = "  +0.58%  "
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "28.01"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.83%  "
# Row 14
$ws.Range("D14").Value = "3.030.57"
$ws.Range("E14").Value = "  +1.45%  "
# Row 15
$ws.Range("D15").Value = "63.085.01"
$ws.Range("E15").Value = "  +0.31%  "
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000147"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.00%  "
# Row 17
$ws.Range("D17").Value = "2.558.63"
$ws.Range("E17").Value = "  +1.00%  "
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.46"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.89%  "
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "341.80"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.90%  "
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.40"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.34%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.85"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.34%  "
# Row 22
$ws.Range("E22").Value = "  -0.02%  "
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.05"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.23%  "
# Row 24
$ws.Range("B24").Value = "Fetch.AI"
$ws.Range("C24").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.65"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.06%  "
# Row 25
$ws.Range("B25").Value = "WrappedeETH"
$ws.Range("C25").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D25").Value = "2.688.78"
$ws.Range("E25").Value = "  +1.29%  "
# Row 26
$ws.Range("E26").Value = "  +1.03%  "
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.22"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +14.04%  "
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.54"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.06%  "
# Row 29
$ws.Range("E29").Value = "  -1.63%  "
# Row 30
$ws.Range("E30").Value = "  +0.39%  "
# Row 31
$ws.Range("E31").Value = "  +6.10%  "
# Row 32
$ws.Range("D32").Value = "0.0₃0832"
$ws.Range("E32").Value = "  +2.31%  "
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "177.55"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.02%  "
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "438.18"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.88%  "
# Row 35
$ws.Range("E35").Value = "  +0.58%  "
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.408"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.86%  "
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "19.35"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.65%  "
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.55"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.36%  "
# Row 39
$ws.Range("E39").Value = "  +0.04%  "
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.76"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.04%  "
# Row 41
$ws.Range("E41").Value = "  +0.13%  "
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "152.33"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.25%  "
# Row 43
$ws.Range("E43").Value = "  +2.07%  "
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "21.37"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.39%  "
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0557"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +6.90%  "
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.609"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.34%  "
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0978"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.09%  "
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0243"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.73%  "
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "18.48"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.59%  "
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.75"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.61%  "
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "11.39"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.11%  "
